$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "countries"
$ws.Name = "countries"

# Update player counts (column E) for specific countries
$ws.Range("E4").Value = 551
$ws.Range("E10").Value = 5225
$ws.Range("E76").Value = 16577
$ws.Range("E105").Value = 534
$ws.Range("E134").Value = 2514
$ws.Range("E214").Value = 7106
